# Update "想去人数" (want-to-go count) values in the F column of the
# "展览" and "全部类型" sheets to reflect the latest scraped counts.

$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 10849
    7  = 175
    9  = 8306
    10 = 42
    12 = 611
    15 = 3308
    17 = 329
    18 = 34
    19 = 797
    21 = 1076
    22 = 288
    23 = 120
    24 = 1797
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
